$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localized rows (14-20) for the "send personal best" feature strings
$data = New-Object 'object[,]' 7,6
$data[0,0] = "facebook friends ranking"
$data[0,1] = "classement amis facebook"
$data[0,2] = "ترتيب أصدقاء الفيسبوك"
$data[0,3] = "Facebook vrienden ranking"
$data[0,4] = "Ranking de amigos de facebook"
$data[0,5] = "classifica amici di facebook"
$data[1,0] = "No Time"
$data[1,1] = "pas le temps"
$data[1,2] = "لا وقت"
$data[1,3] = "geen tijd"
$data[1,4] = "no hay tiempo"
$data[1,5] = "non c'è tempo"
$data[2,0] = "Congratulation You Win with Record of"
$data[2,1] = "Félicitation vous gagnez avec enregistrement de"
$data[2,2] = "مبروك فوزك بسجل"
$data[2,3] = "Gefeliciteerd, je wint met een record van"
$data[2,4] = "Felicidades Ganas con Récord de"
$data[2,5] = "Congratulazioni hai vinto con record di"
$data[3,0] = "all operation are consumed try again"
$data[3,1] = "toutes les opérations sont consommées réessayez"
$data[3,2] = "يتم استهلاك جميع العمليات حاول مرة أخرى"
$data[3,3] = "alle bewerkingen zijn verbruikt probeer het opnieuw"
$data[3,4] = "todas las operaciones se han consumido inténtalo de nuevo"
$data[3,5] = "tutte le operazioni sono state consumate riprovare"
$data[4,0] = "You Lost"
$data[4,1] = "Tu as perdu"
$data[4,2] = "لقد خسرت"
$data[4,3] = "je hebt verloren"
$data[4,4] = "Perdiste"
$data[4,5] = "Hai perso"
$data[5,0] = "Back"
$data[5,1] = "Arrière"
$data[5,2] = "خلف"
$data[5,3] = "Rug"
$data[5,4] = "Atrás"
$data[5,5] = "Di ritorno"
$data[6,0] = "Ops you get decimal value try again"
$data[6,1] = $null
$data[6,2] = "Ops تحصل على قيمة عشرية حاول مرة أخرى"
$data[6,3] = "Oeps, je krijgt een decimale waarde, probeer het opnieuw"
$data[6,4] = "Ops, obtienes un valor decimal, inténtalo de nuevo"
$data[6,5] = "Ops ottieni un valore decimale riprova"

$ws.Range("A14:F20").Value = $data

$ws.Range("A27").Select()
